$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix accessLevel values to be descriptive strings instead of numeric scores,
# and fix the "failclose" -> "failclosed" typo in failMode values.
$ws.Range("E26").Value = '{"accessLevel": "visitor", "availability": 99.0, "failMode": "failclosed"}'
$ws.Range("E27").Value = '{"accessLevel": "staff", "availability": 90.0, "failMode": "failclosed"}'
$ws.Range("E28").Value = '{"accessLevel": "staff", "availability": 98.0,  "failMode": "failopen"}'
$ws.Range("E29").Value = '{"accessLevel": "security", "availability": 80.0, "failMode": "failclosed"}'

# Update the selected/active cell in the sheet view to E31.
$ws.Range("E31").Select()
